$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Georgia's completed-patient count was revised
$ws.Range("B5").Value = 3091

# Ukraine is now complete too: record its patient count and give it the
# same "completed" green fill already used by the other finished rows
# (e.g. Belarus, B3) instead of the old "pending" orange fill.
$ws.Range("B14").Value = 3379
$ws.Range("B14").Interior.Color = $ws.Range("B3").Interior.Color()

# Move the active cell selection
$ws.Range("G14").Select()
